$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Update the two driver assumption cells on Sheet2 ---
# R20: discount/growth rate used for K14:O14 (J38*$R$20 etc.) 0.04 -> 0.02
$ws2.Range("R20").Value = 0.02

# R22: WACC/discount rate used in the NPV formula  8.5% -> 8%
$ws2.Range("R22").Value = 0.08

# --- Re-enter K2's formula so it is no longer grouped as a shared formula
#     together with L2:O2 (those already have their own shared-formula groups) ---
$ws2.Range("K2").Formula = "=J2*1.03"

# --- Selection / active sheet bookkeeping ---
# Sheet1 selection stays where it was (D3); it is simply no longer the tab shown.
$ws1.Range("D3").Select() | Out-Null

# Sheet2's selection moves from S22 to R22, and Sheet2 becomes the active/tab-selected sheet.
$ws2.Activate()
$ws2.Range("R22").Select() | Out-Null
